# Add a "Hospitalized at baseline" derived variable row to the table.
# This inserts a new row at position 61 (pushing the existing rows 61-113
# down to 62-114), fills in the new row's data, resizes the Excel Table
# (ListObject) to include the new row, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new worksheet row above the current row 61, shifting rows down.
$ws.Rows("61:61").Insert()

# Grow the table definition to include the newly inserted row.
$lo.Resize($ws.Range("A1:E114"))

# Populate the new row's cells.
$ws.Range("A61").Value = "O02a"
$ws.Range("B61").Value = "hosp_bl"
$ws.Range("C61").Value = "Outcome"
$ws.Range("D61").Value = "Hospitalized at baseline (within 30 days of diagnosis)"

# Update selection to match the edited cell.
$ws.Range("D61").Select()
